{"js": "// Update patient / mother / report data to the new dataset, and replace\n// the full exam body text with the new placeholder text, per the commit\n// \"Atualizado nova base de dados\".\n\nconst body = context.document.body;\n\nasync function replaceOnce(findText, newText) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + findText);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n\nawait replaceOnce(\"OTAVIO RAMOS DE ALMEIDA   \", \"WILDER VAGNER GASPAR   \");\nawait replaceOnce(\"15/02/1988   \", \"06/10/1984   \");\nawait replaceOnce(\"26294\", \"56276\");\nawait replaceOnce(\"MARIA IRACY RAMOS DOS SANTOS   \", \"MARIA DE FATIMA ARAUJO   \");\nawait replaceOnce(\"20/02/2019   \", \"29/03/2019   \");\n\n// NOTE: deliberately omit the very last trailing \"\\n\" from both the\n// search text and the replacement text. The run's final \"\\n\" character\n// sits at the absolute end of the document story; insertText(\"Replace\")\n// on a range that reaches the story end treats an embedded trailing \"\\n\"\n// in the replacement as a paragraph break instead of a literal newline\n// character. Leaving that last \"\\n\" untouched in the original run keeps\n// it as a literal character in <w:t>, which is what the target XML has\n// (\"Modelo de Laudo.\\n\" as a single run, no extra <w:p>).\nconst oldReport =\n  \".\\n\" +\n  \"T\u00c9CNICA\\n\" +\n  \"Exame  realizado  com cortes tomogr\u00e1ficos computadorizados axiais,  sem a infus\u00e3o endovenosa de contraste iodado, segundo solicita\u00e7\u00e3o do m\u00e9dico assistente.\\n\" +\n  \"Salientamos que a n\u00e3o utiliza\u00e7\u00e3o do meio de contraste iodado por via endovenosa prejudica a adequada caracteriza\u00e7\u00e3o das estruturas abdominais.\\n\" +\n  \"AN\u00c1LISE\\n\" +\n  \"F\u00edgado de topografia, morfologia, situa\u00e7\u00e3o e dimens\u00f5es, preservadas, exibindo coeficientes de atenua\u00e7\u00e3o homog\u00eaneos.\\n\" +\n  \"N\u00e3o h\u00e1 evid\u00eancia de dilata\u00e7\u00e3o das vias biliares intra ou extra-hep\u00e1ticas, bem como da ves\u00edcula biliar.\\n\" +\n  \"Ba\u00e7o, p\u00e2ncreas e adrenais com topografia, dimens\u00f5es, contornos e densidade normais.\\n\" +\n  \"Rins de topografia, morfologia e dimens\u00f5es preservadas, com coeficientes de atenua\u00e7\u00e3o homog\u00eaneos, sem a caracteriza\u00e7\u00e3o de hidronefrose.\\n\" +\n  \"Dois c\u00e1lculos n\u00e3o obstrutivos no ter\u00e7o m\u00e9dio e inferior do rim esquerdo medindo at\u00e9 0,4 cm.\\n\" +\n  \"Aorta e veia cava inferior com posi\u00e7\u00e3o e calibre normais.\\n\" +\n  \"Aus\u00eancia de linfonodomegalias, l\u00edquido livre ou de cole\u00e7\u00f5es organizadas na cavidade abdominal.\\n\" +\n  \"Bexiga urin\u00e1ria em pequena reple\u00e7\u00e3o, com paredes lisas e regulares e conte\u00fado homog\u00eaneo.\\n\" +\n  \"Pr\u00f3stata e ves\u00edculas seminais sem altera\u00e7\u00f5es detect\u00e1veis ao m\u00e9todo.\\n\" +\n  \"OPINI\u00c3O\\n\" +\n  \"Nefrolit\u00edase esquerda n\u00e3o obstrutiva.\\n\" +\n  \"           Dra. Amanda Prist\\n\" +\n  \"             CRM-MG: 56.487\";\n\nconst newReport = \"Modelo de Laudo.\";\n\nawait replaceOnce(oldReport, newReport);\n", "ps1": "# Update patient / mother / report data to the new dataset, and replace\n# the full exam body text with the new placeholder text, per the commit\n# \"Atualizado nova base de dados\".\n\n$d = $word.ActiveDocument\n\n# NOTE: this host's PowerShell interpreter only binds POSITIONAL\n# arguments on custom functions (named `-Param value` args are silently\n# dropped) -- so `Replace-Text` is always called positionally below.\nfunction Replace-Text($FindText, $ReplaceWith) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $FindText\n    $find.Replacement.Text = $ReplaceWith\n\n    $ok = $find.Execute(\n        $FindText,      # FindText\n        $true,          # MatchCase\n        $false,         # MatchWholeWord\n        $false,         # MatchWildcards\n        $false,         # MatchSoundsLike\n        $false,         # MatchAllWordForms\n        $true,          # Forward\n        \"wdFindContinue\", # Wrap\n        $false,         # Format\n        $ReplaceWith,   # ReplaceWith\n        \"wdReplaceAll\"  # Replace\n    )\n\n    if (-not $ok) {\n        throw \"Find/Replace failed for: $FindText\"\n    }\n}\n\nReplace-Text \"OTAVIO RAMOS DE ALMEIDA   \" \"WILDER VAGNER GASPAR   \"\nReplace-Text \"15/02/1988   \" \"06/10/1984   \"\nReplace-Text \"26294\" \"56276\"\nReplace-Text \"MARIA IRACY RAMOS DOS SANTOS   \" \"MARIA DE FATIMA ARAUJO   \"\nReplace-Text \"20/02/2019   \" \"29/03/2019   \"\n\n# The exam report block is the final run of the document, ending at the\n# very end of the story. Its text (below) ends with a literal \"\\n\"\n# character that is NOT part of this Find/Replace span -- deliberately\n# leaving the very last \"\\n\" out of both FindText and ReplaceWith keeps\n# it untouched in the run, exactly as in the target XML\n# (\"Modelo de Laudo.\\n\" as a single run, no extra paragraph inserted).\n$oldReport = @\"\n.\nT\u00c9CNICA\nExame  realizado  com cortes tomogr\u00e1ficos computadorizados axiais,  sem a infus\u00e3o endovenosa de contraste iodado, segundo solicita\u00e7\u00e3o do m\u00e9dico assistente.\nSalientamos que a n\u00e3o utiliza\u00e7\u00e3o do meio de contraste iodado por via endovenosa prejudica a adequada caracteriza\u00e7\u00e3o das estruturas abdominais.\nAN\u00c1LISE\nF\u00edgado de topografia, morfologia, situa\u00e7\u00e3o e dimens\u00f5es, preservadas, exibindo coeficientes de atenua\u00e7\u00e3o homog\u00eaneos.\nN\u00e3o h\u00e1 evid\u00eancia de dilata\u00e7\u00e3o das vias biliares intra ou extra-hep\u00e1ticas, bem como da ves\u00edcula biliar.\nBa\u00e7o, p\u00e2ncreas e adrenais com topografia, dimens\u00f5es, contornos e densidade normais.\nRins de topografia, morfologia e dimens\u00f5es preservadas, com coeficientes de atenua\u00e7\u00e3o homog\u00eaneos, sem a caracteriza\u00e7\u00e3o de hidronefrose.\nDois c\u00e1lculos n\u00e3o obstrutivos no ter\u00e7o m\u00e9dio e inferior do rim esquerdo medindo at\u00e9 0,4 cm.\nAorta e veia cava inferior com posi\u00e7\u00e3o e calibre normais.\nAus\u00eancia de linfonodomegalias, l\u00edquido livre ou de cole\u00e7\u00f5es organizadas na cavidade abdominal.\nBexiga urin\u00e1ria em pequena reple\u00e7\u00e3o, com paredes lisas e regulares e conte\u00fado homog\u00eaneo.\nPr\u00f3stata e ves\u00edculas seminais sem altera\u00e7\u00f5es detect\u00e1veis ao m\u00e9todo.\nOPINI\u00c3O\nNefrolit\u00edase esquerda n\u00e3o obstrutiva.\n           Dra. Amanda Prist\n             CRM-MG: 56.487\n\"@\n\nReplace-Text $oldReport \"Modelo de Laudo.\"\n"}
